$d = $word.ActiveDocument

$replacements = @(
    @("2025-05-31 Saturday", "2025-06-01 Sunday"),
    @("29×19=551", "22×27=594"),
    @("61×60=3660", "59×28=1652"),
    @("68×65=4420", "62×25=1550"),
    @("23×99=2277", "60×28=1680"),
    @("73×89=6497", "64×55=3520"),
    @("45×27=1215", "62×31=1922"),
    @("27×67=1809", "81×34=2754"),
    @("65×99=6435", "82×83=6806"),
    @("30×84=2520", "52×54=2808"),
    @("54×21=1134", "45×24=1080"),
    @("45×60=2700", "46×39=1794"),
    @("84×53=4452", "34×51=1734"),
    @("34×29=986", "63×62=3906"),
    @("34×91=3094", "60×80=4800"),
    @("28×12=336", "48×41=1968"),
    @("57×95=5415", "29×65=1885"),
    @("98×27=2646", "69×49=3381"),
    @("20×89=1780", "95×21=1995"),
    @("85×55=4675", "36×60=2160"),
    @("31×52=1612", "46×41=1886"),
    @("83×73=6059", "80×88=7040"),
    @("60×94=5640", "32×93=2976"),
    @("51×25=1275", "31×64=1984"),
    @("44×20=880", "60×16=960"),
    @("58×22=1276", "42×99=4158")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
